# Generate Report for handoff
#
# Re-generates the localization-status report: the file
# d71034b0-0ab0-4cb6-a7cb-662986bb4024 moves from "Ready for handoff" into
# "In Translation" (its previous handoff timestamp is retained), while
# 5aee7500-ca76-4eca-adc0-d8a683b8b5e0 is freshly re-handed-off (status stays
# "Ready for handoff" but gets a new, later "Latest Handoff Datetime"). The
# report rows are re-emitted with the d71034b0 entry first (the hyperlink
# targets stay anchored to their original cell position - only the visible
# link text is refreshed to match the regenerated row).

$wb = $excel.ActiveWorkbook

$repoBase   = "https://github.com/OpenLocalizationTest/oltest/blob/cb06ec6f780c57fbb46fcd4726d210ab3b595678"
$zhHandoff  = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/ca954d6be9c8d67c0e421d3ba96ba0ff53e54312/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/xinjiang/ht"
$deHandoff  = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/50d0e3579a7c5f2dedaefdd70fed79598b1ff6c3/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/xinjiang/ht"

$d71 = "d71034b0-0ab0-4cb6-a7cb-662986bb4024"
$ae5 = "5aee7500-ca76-4eca-adc0-d8a683b8b5e0"

$d71Md  = "$d71.md"
$ae5Md  = "$ae5.md"

$d71ZhXlf = "$d71.394dcf54d44ec9b8eb57c7af0afa8a8aa1f7976d.zh-cn.xlf"
$d71DeXlf = "$d71.394dcf54d44ec9b8eb57c7af0afa8a8aa1f7976d.de-de.xlf"
$ae5ZhXlf = "$ae5.f2215be77fa5769f6285a39f8c93f997f9a6a744.zh-cn.xlf"
$ae5DeXlf = "$ae5.f2215be77fa5769f6285a39f8c93f997f9a6a744.de-de.xlf"

# ---------------------------------------------------------------------------
# Overview sheet: d71034b0 row now listed first and flips to "In Translation"
# ---------------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")

$wsOverview.Range("A2").Value = $d71Md
$wsOverview.Range("B2").Value = "In Translation"
$wsOverview.Range("C2").Value = "In Translation"

$wsOverview.Range("A3").Value = $ae5Md
$wsOverview.Range("B3").Value = "Ready for handoff"
$wsOverview.Range("C3").Value = "Ready for handoff"

# Hyperlink relationship targets stay bound to the cell position (the
# underlying .rels entries are untouched); only the rendered text changes.
$wsOverview.Hyperlinks.Delete()
$wsOverview.Hyperlinks.Add($wsOverview.Range("A2"), "$repoBase/e2e/$ae5Md", $null, $null, $d71Md)
$wsOverview.Hyperlinks.Add($wsOverview.Range("A3"), "$repoBase/e2e/$d71Md", $null, $null, $ae5Md)
$wsOverview.Hyperlinks.Add($wsOverview.Range("A4"), "$repoBase/.localization-config", $null, $null, ".localization-config")

# ---------------------------------------------------------------------------
# zh-cn sheet
# ---------------------------------------------------------------------------
$wsZh = $wb.Worksheets.Item("zh-cn")

$wsZh.Range("A2").Value = $d71Md
$wsZh.Range("B2").Value = "In Translation"
$wsZh.Range("C2").Value = $d71ZhXlf
$wsZh.Range("D2").Value = "2016-01-26 11:23:44"

$wsZh.Range("A3").Value = $ae5Md
$wsZh.Range("B3").Value = "Ready for handoff"
$wsZh.Range("C3").Value = $ae5ZhXlf
$wsZh.Range("D3").Value = "2016-01-26 11:25:20"

$wsZh.Hyperlinks.Delete()
$wsZh.Hyperlinks.Add($wsZh.Range("A2"), "$repoBase/e2e/$ae5Md", $null, $null, $d71Md)
$wsZh.Hyperlinks.Add($wsZh.Range("C2"), "$zhHandoff/$ae5ZhXlf", $null, $null, $d71ZhXlf)
$wsZh.Hyperlinks.Add($wsZh.Range("A3"), "$repoBase/e2e/$d71Md", $null, $null, $ae5Md)
$wsZh.Hyperlinks.Add($wsZh.Range("C3"), "$zhHandoff/$d71ZhXlf", $null, $null, $ae5ZhXlf)
$wsZh.Hyperlinks.Add($wsZh.Range("A4"), "$repoBase/.localization-config", $null, $null, ".localization-config")

# ---------------------------------------------------------------------------
# de-de sheet
# ---------------------------------------------------------------------------
$wsDe = $wb.Worksheets.Item("de-de")

$wsDe.Range("A2").Value = $d71Md
$wsDe.Range("B2").Value = "In Translation"
$wsDe.Range("C2").Value = $d71DeXlf
$wsDe.Range("D2").Value = "2016-01-26 11:23:56"

$wsDe.Range("A3").Value = $ae5Md
$wsDe.Range("B3").Value = "Ready for handoff"
$wsDe.Range("C3").Value = $ae5DeXlf
$wsDe.Range("D3").Value = "2016-01-26 11:25:30"

$wsDe.Hyperlinks.Delete()
$wsDe.Hyperlinks.Add($wsDe.Range("A2"), "$repoBase/e2e/$ae5Md", $null, $null, $d71Md)
$wsDe.Hyperlinks.Add($wsDe.Range("C2"), "$deHandoff/$ae5DeXlf", $null, $null, $d71DeXlf)
$wsDe.Hyperlinks.Add($wsDe.Range("A3"), "$repoBase/e2e/$d71Md", $null, $null, $ae5Md)
$wsDe.Hyperlinks.Add($wsDe.Range("C3"), "$deHandoff/$d71DeXlf", $null, $null, $ae5DeXlf)
$wsDe.Hyperlinks.Add($wsDe.Range("A4"), "$repoBase/.localization-config", $null, $null, ".localization-config")
